$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the survey question text (line-break moved earlier in the sentence)
$ws.Range("A2").Value = "`"Governments should actively cooperate to have all countries`nconverge in terms of GDP per capita by the end of the century`""

# Update the recalculated convergence support shares ("final data")
$ws.Range("B2").Value = 0.704965329416964
$ws.Range("K2").Value = 0.703874631903231
$ws.Range("L2").Value = 0.777916745185535
$ws.Range("N2").Value = 0.562406199574745
